$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 4692
$ws.Range("J3").Value = 4970
$ws.Range("J4").Value = 1103
$ws.Range("J5").Value = 393
$ws.Range("J6").Value = 6066
$ws.Range("J7").Value = 17224
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J3").Value = 6
$ws.Range("J7").Value = 498
$ws.Range("J8").Value = 1109
$ws.Range("J9").Value = 88
$ws.Range("J10").Value = 108
$ws.Range("J11").Value = 261
$ws.Range("J12").Value = 36
$ws.Range("J15").Value = 189
$ws.Range("J19").Value = 502
$ws.Range("J20").Value = 360
$ws.Range("J29").Value = 984
$ws.Range("J31").Value = 148
$ws.Range("J33").Value = 782
$ws.Range("J34").Value = 81
$ws.Range("J36").Value = 238
$ws.Range("J37").Value = 542
$ws.Range("J38").Value = 9
$ws.Range("J43").Value = 144
$ws.Range("J44").Value = 126
$ws.Range("J48").Value = 187
$ws.Range("J50").Value = 95
$ws.Range("J51").Value = 219
$ws.Range("J53").Value = 214
$ws.Range("J56").Value = 23
$ws.Range("J60").Value = 111
$ws.Range("J63").Value = 65
$ws.Range("J65").Value = 461
$ws.Range("J66").Value = 52
$ws.Range("J67").Value = 666
$ws.Range("J68").Value = 30
$ws.Range("J71").Value = 56
$ws.Range("J73").Value = 156
$ws.Range("J78").Value = 217
$ws.Range("J79").Value = 494
$ws.Range("J83").Value = 377
$ws.Range("J84").Value = 145
$ws.Range("J85").Value = 766
$ws.Range("J90").Value = 199
$ws.Range("J97").Value = 136
$ws.Range("J99").Value = 256
$ws.Range("J101").Value = 17224
$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J3").Value = 152
$ws.Range("J6").Value = 157
$ws.Range("J7").Value = 498
$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("J4").Value = 25
$ws.Range("J6").Value = 65
$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("J6").Value = 95
$ws.Range("J7").Value = 261
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J2").Value = 317
$ws.Range("J3").Value = 336
$ws.Range("J6").Value = 365
$ws.Range("J7").Value = 1109
$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("J2").Value = 45
$ws.Range("J3").Value = 33
$ws.Range("J6").Value = 128
$ws.Range("J7").Value = 214
$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J3").Value = 277
$ws.Range("J6").Value = 217
$ws.Range("J7").Value = 766
$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J2").Value = 113
$ws.Range("J3").Value = 140
$ws.Range("J7").Value = 377
$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("J2").Value = 71
$ws.Range("J6").Value = 70
$ws.Range("J7").Value = 256
$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J2").Value = 161
$ws.Range("J3").Value = 263
$ws.Range("J7").Value = 666
$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("J6").Value = 38
$ws.Range("J7").Value = 148
$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("J3").Value = 45
$ws.Range("J7").Value = 145
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J3").Value = 189
$ws.Range("J6").Value = 157
$ws.Range("J7").Value = 542
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J2").Value = 194
$ws.Range("J3").Value = 253
$ws.Range("J7").Value = 782
$ws = $wb.Worksheets.Item("New City")
$ws.Range("J3").Value = 139
$ws.Range("J6").Value = 161
$ws.Range("J7").Value = 461
$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J2").Value = 298
$ws.Range("J3").Value = 339
$ws.Range("J4").Value = 56
$ws.Range("J6").Value = 251
$ws.Range("J7").Value = 984
$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J3").Value = 144
$ws.Range("J6").Value = 185
$ws.Range("J7").Value = 502
$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("J6").Value = 41
$ws.Range("J7").Value = 126
$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("J2").Value = 29
$ws.Range("J7").Value = 187
$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("J4").Value = 4
$ws.Range("J7").Value = 108
$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("J6").Value = 57
$ws.Range("J7").Value = 217
$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J2").Value = 141
$ws.Range("J6").Value = 135
$ws.Range("J7").Value = 494
$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("J4").Value = 10
$ws.Range("J6").Value = 40
$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J4").Value = 34
$ws.Range("J6").Value = 95
$ws.Range("J7").Value = 360
$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("J3").Value = 74
$ws.Range("J7").Value = 238
$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("J2").Value = 24
$ws.Range("J7").Value = 81
$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("J4").Value = 7
$ws.Range("J7").Value = 189
$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("J2").Value = 26
$ws.Range("J7").Value = 95
$ws = $wb.Worksheets.Item("North Center")
$ws.Range("J2").Value = 10
$ws.Range("J7").Value = 52
$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("J2").Value = 24
$ws.Range("J6").Value = 32
$ws.Range("J7").Value = 88
$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("J2").Value = 55
$ws.Range("J6").Value = 44
$ws.Range("J7").Value = 156
$ws = $wb.Worksheets.Item("West Town")
$ws.Range("J3").Value = 16
$ws.Range("J7").Value = 136
$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("J5").Value = 5
$ws.Range("J7").Value = 199
$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("J2").Value = 52
$ws.Range("J7").Value = 219
$ws = $wb.Worksheets.Item("North Park")
$ws.Range("J6").Value = 7
$ws.Range("J7").Value = 30
$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("J2").Value = 38
$ws.Range("J7").Value = 111
$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("J4").Value = 16
$ws.Range("J7").Value = 144
$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("J2").Value = 14
$ws.Range("J7").Value = 56
$ws = $wb.Worksheets.Item("Magnificent Mile")
$ws.Range("J6").Value = 13
$ws.Range("J7").Value = 23
$ws = $wb.Worksheets.Item("Andersonville")
$ws.Range("J3").Value = 2
$ws.Range("J6").Value = 6
$ws = $wb.Worksheets.Item("Beverly")
$ws.Range("J6").Value = 24
$ws.Range("J7").Value = 36
$ws = $wb.Worksheets.Item("Grant Park")
$ws.Range("J2").Value = 5
$ws.Range("J6").Value = 9
